$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5", "D6", "D9", "D11", "D12", "D13", "D14", "D16", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D45", "D47", "D48", "D49", "D51")
foreach ($cellAddr in $textCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

$ws.Range("D2").Value = "62.837.93"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "3.140.70"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "595.41"
$ws.Range("E5").Value = "  +1.58%  "
$ws.Range("D6").Value = "134.05"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "3.138.49"
$ws.Range("E8").Value = "  -1.25%  "
$ws.Range("D9").Value = "0.507"
$ws.Range("E9").Value = "  +1.25%  "
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("D11").Value = "5.35"
$ws.Range("E11").Value = "  +2.69%  "
$ws.Range("D12").Value = "0.449"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").Value = "0.0000241"
$ws.Range("E13").Value = "  +3.12%  "
$ws.Range("D14").Value = "34.38"
$ws.Range("E14").Value = "  +4.11%  "
$ws.Range("D15").Value = "3.664.24"
$ws.Range("E15").Value = "  -1.05%  "
$ws.Range("D16").Value = "0.120"
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("D17").Value = "3.149.94"
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("D18").Value = "62.951.17"
$ws.Range("E18").Value = "  +0.86%  "
$ws.Range("D19").Value = "6.50"
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("D20").Value = "458.07"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("D21").Value = "13.83"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "0.689"
$ws.Range("E22").Value = "  -1.51%  "
$ws.Range("D23").Value = "7.55"
$ws.Range("E23").Value = "  -0.62%  "
$ws.Range("D24").Value = "13.14"
$ws.Range("E24").Value = "  -1.76%  "
$ws.Range("D25").Value = "82.52"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("B27").Value = "FirstDigitalUSD"
$ws.Range("C27").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "2.66"
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("D29").Value = "2.06"
$ws.Range("E29").Value = "  +3.06%  "
$ws.Range("D30").Value = "7.61"
$ws.Range("E30").Value = "  -2.14%  "
$ws.Range("D31").Value = "6.62"
$ws.Range("E31").Value = "  -4.40%  "
$ws.Range("D32").Value = "26.88"
$ws.Range("E32").Value = "  -1.20%  "
$ws.Range("D33").Value = "0.0997"
$ws.Range("E33").Value = "  -1.66%  "
$ws.Range("D34").Value = "2.38"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("D36").Value = "5.84"
$ws.Range("E36").Value = "  +1.29%  "
$ws.Range("D37").Value = "50.94"
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").Value = "0.0₃0719"
$ws.Range("E38").Value = "  +4.33%  "
$ws.Range("D39").Value = "0.0386"
$ws.Range("E39").Value = "  +0.96%  "
$ws.Range("D40").Value = "8.08"
$ws.Range("E40").Value = "  +1.49%  "
$ws.Range("D41").Value = "0.111"
$ws.Range("E41").Value = "  -0.48%  "
$ws.Range("D42").Value = "2.59"
$ws.Range("E42").Value = "  -1.29%  "
$ws.Range("D43").Value = "386.49"
$ws.Range("E43").Value = "  -5.71%  "
$ws.Range("D44").Value = "2.765.34"
$ws.Range("E44").Value = "  -5.63%  "
$ws.Range("D45").Value = "0.247"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "125.81"
$ws.Range("E47").Value = "  +1.04%  "
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").Value = "35.38"
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("B49").Value = "Fetch.AI"
$ws.Range("C49").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D49").Value = "2.09"
$ws.Range("E49").Value = "  -2.03%  "
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("D51").Value = "24.65"
$ws.Range("E51").Value = "  -2.50%  "

foreach ($cellAddr in $textCells) {
    $ws.Range($cellAddr).Style = "Normal"
}
